$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, reflecting repulled data / recalculated mean.
$updates = @{
    2  = 1
    4  = -2
    5  = -1
    6  = -5
    7  = -3
    8  = -5
    9  = -6
    10 = -3
    11 = 3
    12 = -4
    13 = -5
    14 = -5
    15 = 1
    17 = -2
    18 = -3
    19 = 2
    20 = 7
    21 = -2
    22 = 2
    23 = -3
    24 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
